$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from G1 (existing header cell) to H1, then set the value.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
